# Update the two-digit multiplication answer cells in the table
# with the newly generated values (commit c986bee).
$d = $word.ActiveDocument

$d.Content.Find.Execute("32×86=2752", $true, $false, $false, $false, $false, $true, 1, $false, "60×60=3600", 2) | Out-Null
$d.Content.Find.Execute("81×30=2430", $true, $false, $false, $false, $false, $true, 1, $false, "88×23=2024", 2) | Out-Null
$d.Content.Find.Execute("73×55=4015", $true, $false, $false, $false, $false, $true, 1, $false, "76×29=2204", 2) | Out-Null
$d.Content.Find.Execute("90×89=8010", $true, $false, $false, $false, $false, $true, 1, $false, "35×33=1155", 2) | Out-Null
$d.Content.Find.Execute("27×83=2241", $true, $false, $false, $false, $false, $true, 1, $false, "71×92=6532", 2) | Out-Null
$d.Content.Find.Execute("97×59=5723", $true, $false, $false, $false, $false, $true, 1, $false, "84×95=7980", 2) | Out-Null
$d.Content.Find.Execute("63×22=1386", $true, $false, $false, $false, $false, $true, 1, $false, "27×23=621", 2) | Out-Null
$d.Content.Find.Execute("94×52=4888", $true, $false, $false, $false, $false, $true, 1, $false, "74×26=1924", 2) | Out-Null
$d.Content.Find.Execute("37×26=962", $true, $false, $false, $false, $false, $true, 1, $false, "17×92=1564", 2) | Out-Null
$d.Content.Find.Execute("57×62=3534", $true, $false, $false, $false, $false, $true, 1, $false, "87×47=4089", 2) | Out-Null
$d.Content.Find.Execute("64×53=3392", $true, $false, $false, $false, $false, $true, 1, $false, "88×17=1496", 2) | Out-Null
$d.Content.Find.Execute("74×52=3848", $true, $false, $false, $false, $false, $true, 1, $false, "41×52=2132", 2) | Out-Null
$d.Content.Find.Execute("30×96=2880", $true, $false, $false, $false, $false, $true, 1, $false, "80×78=6240", 2) | Out-Null
$d.Content.Find.Execute("54×52=2808", $true, $false, $false, $false, $false, $true, 1, $false, "67×23=1541", 2) | Out-Null
$d.Content.Find.Execute("54×47=2538", $true, $false, $false, $false, $false, $true, 1, $false, "73×88=6424", 2) | Out-Null
$d.Content.Find.Execute("64×23=1472", $true, $false, $false, $false, $false, $true, 1, $false, "38×95=3610", 2) | Out-Null
$d.Content.Find.Execute("37×87=3219", $true, $false, $false, $false, $false, $true, 1, $false, "32×39=1248", 2) | Out-Null
$d.Content.Find.Execute("99×23=2277", $true, $false, $false, $false, $false, $true, 1, $false, "25×44=1100", 2) | Out-Null
$d.Content.Find.Execute("14×28=392", $true, $false, $false, $false, $false, $true, 1, $false, "23×19=437", 2) | Out-Null
$d.Content.Find.Execute("19×15=285", $true, $false, $false, $false, $false, $true, 1, $false, "94×51=4794", 2) | Out-Null
$d.Content.Find.Execute("77×93=7161", $true, $false, $false, $false, $false, $true, 1, $false, "46×74=3404", 2) | Out-Null
$d.Content.Find.Execute("76×91=6916", $true, $false, $false, $false, $false, $true, 1, $false, "61×64=3904", 2) | Out-Null
$d.Content.Find.Execute("82×42=3444", $true, $false, $false, $false, $false, $true, 1, $false, "28×63=1764", 2) | Out-Null
$d.Content.Find.Execute("20×26=520", $true, $false, $false, $false, $false, $true, 1, $false, "38×59=2242", 2) | Out-Null
$d.Content.Find.Execute("40×62=2480", $true, $false, $false, $false, $false, $true, 1, $false, "28×17=476", 2) | Out-Null
